# Updates cryptos list values (prices & 1h volume change %) per the
# upstream GitHub Actions scrape-and-commit job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.048.22"
$ws.Range("E2").Value = "  -0.13%  "
$ws.Range("D3").Value = "1.817.88"
$ws.Range("E3").Value = "  -0.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9988"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "240.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.01%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6140"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9996"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07311"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2886"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07651"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "1.828.76"
$ws.Range("E12").Value = "  -0.25%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.933"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6578"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "81.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.66%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008983"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.96%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.824"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.58%  "
$ws.Range("D18").Value = "29.049.32"
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").Value = "2.068.14"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "236.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9994"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.096"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.9994"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.63"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1402"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.59%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.399"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.12%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.56"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.74%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.475"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.71%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05546"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.85%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.079"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.085"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.54%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.205"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7326"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.15%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.805"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.92%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.126"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.19%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.612"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.822"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.28%  "
$ws.Range("D39").Value = "1.206.04"
$ws.Range("E39").Value = "  -1.18%  "
$ws.Range("E40").Value = "  -1.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.355"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.69%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8900"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9996"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.77"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").Value = "1.981.49"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "64.31"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000121"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.5079"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05748"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.24%  "

# Row 49/50: EnergySwap and TheSandbox swapped ranking positions.
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.039"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.35%  "

$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.3979"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.31%  "
